$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 39 (pushes existing rows 39..117 down to 40..118,
# extending the table from A1:R117 to A1:R118).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 44665
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 100112052
$ws.Range("G39").Value = "Albahaca"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = 6000
$ws.Range("N39").Value = "$/docena"
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 5000
$ws.Range("Q39").Value = 1.2
$ws.Range("R39").Value = "Hortaliza"
